# Fruta / hortaliza, semanal
# Insert a new daily price record at row 330 (Vega Monumental Concepción -
# Naranja, Valencia, Primera) shifting the existing rows 330..431 down to
# 331..432.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 330; everything from 330 downward shifts
# down by one (old row 330 becomes row 331, ..., old row 431 becomes 432).
$ws.Rows.Item(330).Insert()

# Populate the newly inserted row 330 with the new record's data.
$ws.Cells.Item(330, 1).Value = 11
$ws.Cells.Item(330, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(330, 3).Value = "Bíobío"
$ws.Cells.Item(330, 4).Value = 44985
$ws.Cells.Item(330, 5).Value = 8
$ws.Cells.Item(330, 6).Value = "Fruta"
$ws.Cells.Item(330, 7).Value = 100102
$ws.Cells.Item(330, 8).Value = "Cítricos"
$ws.Cells.Item(330, 9).Value = 100102005
$ws.Cells.Item(330, 10).Value = "Naranja"
$ws.Cells.Item(330, 11).Value = "Valencia"
$ws.Cells.Item(330, 12).Value = "Primera"
$ws.Cells.Item(330, 13).Value = 270
$ws.Cells.Item(330, 14).Value = 11000
$ws.Cells.Item(330, 15).Value = 12000
$ws.Cells.Item(330, 16).Value = 11556
$ws.Cells.Item(330, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(330, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(330, 19).Value = 770
$ws.Cells.Item(330, 20).Value = 15
